$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original "Method" values (currently in column A, rows 2-5)
# before we overwrite them, so we can move them into the new column C.
# (Read via Value2 - the plain Value getter is unreliable for reads in
# this runtime.)
$origMethod2 = $ws.Range("A2").Value2
$origMethod3 = $ws.Range("A3").Value2
$origMethod4 = $ws.Range("A4").Value2
$origMethod5 = $ws.Range("A5").Value2

# New header row: A1=Network, B1=Alpha, C1=Method, D1=Average Significant Percentage
$ws.Range("A1").Value = "Network"
$ws.Range("B1").Value = "Alpha"
$ws.Range("C1").Value = "Method"
$ws.Range("D1").Value = "Average Significant Percentage"

# Apply the same header formatting (bold font, centered, thin border) used
# by the pre-existing header cell to all four header cells, A1:D1, by
# copying its format.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the Method names (moved to column C) explicitly from the values
# captured above, to be safe.
$ws.Range("C2").Value = $origMethod2
$ws.Range("C3").Value = $origMethod3
$ws.Range("C4").Value = $origMethod4
$ws.Range("C5").Value = $origMethod5

# New "Network" column values (A2:A5) - all H_sapiens for this report.
$ws.Range("A2").Value = "H_sapiens"
$ws.Range("A3").Value = "H_sapiens"
$ws.Range("A4").Value = "H_sapiens"
$ws.Range("A5").Value = "H_sapiens"

# New "Alpha" column values (B2:B5).
$ws.Range("B2").Value = 0.1
$ws.Range("B3").Value = 0.1
$ws.Range("B4").Value = 0.1
$ws.Range("B5").Value = 0.1

# Updated "Average Significant Percentage" values (now column D).
$ws.Range("D2").Value = 0.1362397820163488
$ws.Range("D3").Value = 0.04541326067211626
$ws.Range("D4").Value = 0.04541326067211626
$ws.Range("D5").Value = 0.01135331516802906
